$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value while forcing it to remain Text
# (leading apostrophe = Excel "treat as text" quote-prefix marker).

$ws.Range("D2").Value = "'248.99"
$ws.Range("E2").Value = "1BNBBNBBestin24h"

$ws.Range("D3").Value = "'22.10"

$ws.Range("D4").Value = "'5.529"

$ws.Range("D5").Value = "'0.05624"

$ws.Range("D6").Value = "'6.467"

$ws.Range("D8").Value = "'1.041"

$ws.Range("D9").Value = "'0.1432"

$ws.Range("D10").Value = "'0.07323"

$ws.Range("D11").Value = "'0.03115"

$ws.Range("D12").Value = "'0.02917"

$ws.Range("D13").Value = "'0.09269"

$ws.Range("D14").Value = "'0.001668"

$ws.Range("D15").Value = "'3.228"

$ws.Range("D16").Value = "'0.04728"

$ws.Range("D17").Value = "'0.0005814"
$ws.Range("E17").Value = "16OneONE"

$ws.Range("D18").Value = "'0.006388"

$ws.Range("D19").Value = "'0.005068"

$ws.Range("D20").Value = "'0.001054"

$ws.Range("D21").Value = "'0.0001502"

$ws.Range("D22").Value = "'3.977"

$ws.Range("D24").Value = "'2.111"

$ws.Range("D27").Value = "'0.0002986"

$ws.Range("D40").Value = "'0.04149"

# Rows 41-43: coins rotated (Kick -> BKEX -> CEJI -> Kick)
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1042"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002974"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.007069"
$ws.Range("E43").Value = "42KickTokenKICK"

$ws.Range("D44").Value = "'0.008701"

$ws.Range("D45").Value = "'0.00005641"

$ws.Range("D46").Value = "'0.00000000751"

$ws.Range("D47").Value = "'0.6805"

$ws.Range("D48").Value = "'0.01635"
